# Commit: add a third worksheet "sode" (soda / mixers) next to the existing
# "hva" and "alko" sheets, fill "alko" (alcohol) and "sode" with drink
# ingredient lists, bold their headers, and move the active tab from "hva"
# to "alko".

$wb = $excel.ActiveWorkbook

$wsHva  = $wb.Worksheets.Item(1)
$wsAlko = $wb.Worksheets.Item(2)

# Create the new "sode" sheet by copying the (still empty) "alko" sheet, so
# it inherits the same sheet defaults (row height etc.) instead of the
# engine's brand-new-sheet defaults, then rename + clear it out.
$wsAlko.Copy([System.Reflection.Missing]::Value, $wsAlko)
$wsSode = $wb.Worksheets.Item(3)
$wsSode.Name = "sode"

# --- "alko" sheet: fill in alcohol list -----------------------------------
# Typing order matters for shared-string ordering: Vodka/Gin first, then
# going back to fill in the "Alkohol" header.
$wsAlko.Range("A2").Value = "Vodka"
$wsAlko.Range("A3").Value = "Gin"
$wsAlko.Range("A1").Value = "Alkohol"
$wsAlko.Range("A1").Font.Bold = $true

$wsSode.Range("A1").Value = "Opblanding"
$wsSode.Range("A1").Font.Bold = $true
$wsSode.Range("A2").Value = "Applesinjuice"
$wsSode.Range("A3").Value = "Kakao"
$wsSode.Range("A4").Value = "Coca Cola"
$wsSode.Range("A5").Value = "Gazoz"

# --- back to "alko": finish the alcohol list -------------------------------
$wsAlko.Range("A4").Value = "Rom"
$wsAlko.Range("A5").Value = "Hvidrom"
$wsAlko.Range("A6").Value = "Små Sure"

# --- back to "sode": finish the mixer list ---------------------------------
$wsSode.Range("A6").Value = "Fanta"
$wsSode.Range("A7").Value = "Tonic"

# --- column widths -----------------------------------------------------------
# The runtime's ColumnWidth setter quantizes to 1/6-character steps (xml
# width = ColumnWidth + 5/6, rounded to the nearest 1/6), so these land on
# the closest representable width to the authored xlsx (15.36328125,
# 13.453125, 14.453125 respectively).
$wsAlko.Columns.Item(1).ColumnWidth = 14.5
$wsAlko.Columns.Item(2).ColumnWidth = 12.666666666666666
$wsSode.Columns.Item(1).ColumnWidth = 13.666666666666666

# --- page setup on "sode" ---------------------------------------------------
$wsSode.PageSetup.PaperSize = 9
$wsSode.PageSetup.Orientation = 1

# --- selections --------------------------------------------------------------
[void]$wsHva.Range("E1:E1048576").Select()
[void]$wsAlko.Range("D26").Select()
[void]$wsAlko.Select()
